$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value "Hi" in B1 (next to existing "Hey" in A1)
$ws.Range("B1").Value = "Hi"

# Select B1 so that the saved view matches the active cell/selection
$ws.Range("B1").Select()
